$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Incomes")

# The "Table7" structured table holds the incomes data (A:C). Adding a
# ListRow expands the table (and its autofilter) by one row and keeps
# everything else - styles, headers, sort state - intact.
$tbl = $ws.ListObjects.Item("Table7")
$newRow = $tbl.ListRows.Add()

$lastDataRow = $tbl.Range.Rows.Count + $tbl.Range.Row - 2
$newDataRow = $lastDataRow + 1

# Copy number formatting (date / currency) from the previous row so the
# new cells reuse the existing cell styles instead of creating new ones.
$ws.Range("A" + $lastDataRow).Copy() | Out-Null
$ws.Range("A" + $newDataRow).PasteSpecial(-4122) | Out-Null
$ws.Range("B" + $lastDataRow).Copy() | Out-Null
$ws.Range("B" + $newDataRow).PasteSpecial(-4122) | Out-Null

# New income entry: a missed "Parents transfer" row.
$ws.Range("A" + $newDataRow).Value = 45715
$ws.Range("B" + $newDataRow).Value = 150
$ws.Range("C" + $newDataRow).Value = "Parents transfer"

$ws.Range("B89").Select() | Out-Null
